$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are text (e.g. "37.352.98", "15.70") that must
# stay text -- Excel auto-converts plain-number-looking strings assigned
# via .Value into floating point numbers, so force text format first and
# restore the default "Normal" style afterwards to avoid leaving a stray
# number-format style on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.352.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.008.14"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "258.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.51%  "

$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.14%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0768"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.10%  "

$ws.Range("E11").Value = "  -2.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.303.66"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.47%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.21"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.35%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.799"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.88"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.91%  "

$ws.Range("E16").Value = "  -6.02%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.988.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.37%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.279.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.30%  "

$ws.Range("E20").Value = "  -4.19%  "

$ws.Range("E21").Value = "  -2.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.23"
$ws.Range("D22").Style = "Normal"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.98%  "

$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  -1.24%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.84"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.83%  "

$ws.Range("E29").Value = "  -6.28%  "

$ws.Range("E30").Value = "  -5.91%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.119"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("E32").Value = "  -4.59%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.87%  "

$ws.Range("E34").Value = "  +0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.36"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.99%  "

$ws.Range("E36").Value = "  +0.39%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.25%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.39%  "

$ws.Range("E40").Value = "  +4.13%  "

$ws.Range("E41").Value = "  -0.16%  "

$ws.Range("E42").Value = "  -1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0924"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.419.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.08"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.65%  "

$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.19%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -7.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.195.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -8.25%  "
